$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached "last saved" date shown by the auto-updating
#    date placeholders (notes master, slide master, and every slide
#    layout's footer) from 9/16/2017 to 9/17/2017.
# ---------------------------------------------------------------------
function Update-DateShape($shape) {
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "9/16/2017") {
            $tr.Text = "9/17/2017"
        }
    }
}

$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    Update-DateShape $nm.Shapes.Item($i)
}

$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    Update-DateShape $m.Shapes.Item($i)
}

$layouts = $m.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $cl = $layouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        Update-DateShape $cl.Shapes.Item($i)
    }
}

# ---------------------------------------------------------------------
# 2) Slide 16 ("Next Steps"): extend the Guided Practice bullet to cover
#    both 6.6 and 6.7.
# ---------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$content = $s16.Shapes.Item(2)
$tr16 = $content.TextFrame.TextRange
$para = $tr16.Paragraphs(3, 1)

# Keep "Do " as-is and replace "Guided Practice 6.6" with the updated
# wording, which naturally splits the paragraph into the two runs seen
# in the authored edit.
$rest = $para.Characters(4, 19)
$rest.Text = "Guided Practices 6.6 and 6.7"
